# Lecture partielle de l'EDT M1 MIAGE.
# Shift the schedule dates forward by 1096 days (3 years) and rename the
# weekday labels accordingly: "jeudi" -> "lundi" and "dimanche" -> "jeudi".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "A" column holds the week-start date and whose "B" column
# holds the week-day label shared string.
$mondayRows = @(2, 8, 14, 20, 26, 32)
$thursdayRows = @(5, 11, 17, 23, 29, 35)

foreach ($r in $mondayRows) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 + 1096
    $ws.Cells.Item($r, 2).Value = "lundi"
}

foreach ($r in $thursdayRows) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 + 1096
    $ws.Cells.Item($r, 2).Value = "jeudi"
}
